# Insert a new weekly record as row 5, pushing all existing rows (5..113)
# down by one (to 6..114). This matches the observed diff where every row's
# data equals the previous row's data, shifted down by one row, with a brand
# new record occupying the new row 5 and the dimension growing to A1:R114.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 5, shifting rows 5:113 down to 6:114.
$ws.Rows("5:5").Insert()

# Populate the newly inserted row 5 with the new weekly price record.
$ws.Range("A5").Value = 6
$ws.Range("B5").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C5").Value = "Metropolitana"
$ws.Range("D5").Value = 44496
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = 100112001
$ws.Range("G5").Value = "Berenjena"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 1100
$ws.Range("K5").Value = 6000
$ws.Range("L5").Value = 7000
$ws.Range("M5").Value = 6591
$ws.Range("N5").Value = "`$/caja 50 unidades"
$ws.Range("O5").Value = "Región de Arica y Parinacota"
$ws.Range("P5").Value = 132
$ws.Range("Q5").Value = 50
$ws.Range("R5").Value = "Hortaliza"
